$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# --- New headers for columns C (assetAddress) and D (title) ---
$ws.Cells.Item(1,3).Value = "assetAddress"
$ws.Cells.Item(1,4).Value = "title"

# --- Row 2: Item 1 ---
$ws.Cells.Item(2,4).Value = $ws.Cells.Item(2,2).Value()
$ws.Cells.Item(2,2).Value = "_10000_Item_1"
$ws.Cells.Item(2,3).Value = "AssetAddress.CubeBlue_Pickupable"

# --- Row 3: Item 2 ---
$ws.Cells.Item(3,4).Value = $ws.Cells.Item(3,2).Value()
$ws.Cells.Item(3,2).Value = "_10001_Item_2"
$ws.Cells.Item(3,3).Value = "AssetAddress.CubeRed_Pickupable"

# --- Row 4: Puzzle Piece A ---
$ws.Cells.Item(4,4).Value = $ws.Cells.Item(4,2).Value()
$ws.Cells.Item(4,2).Value = "_10100_PuzzleBlock_A"
$ws.Cells.Item(4,3).Value = "AssetAddress.PuzzleBlock_A_Pickupable"

# --- Row 5: Puzzle Piece B ---
$ws.Cells.Item(5,4).Value = $ws.Cells.Item(5,2).Value()
$ws.Cells.Item(5,2).Value = "_10101_PuzzleBlock_B"
$ws.Cells.Item(5,3).Value = "AssetAddress.PuzzleBlock_B_Pickupable"

# --- Row 6: Old Key ---
$ws.Cells.Item(6,4).Value = $ws.Cells.Item(6,2).Value()
$ws.Cells.Item(6,2).Value = "_10201_Key_A"
$ws.Cells.Item(6,3).Value = "AssetAddress.Key_A_Pickupable"

# --- Row 7: Book A ---
$ws.Cells.Item(7,4).Value = $ws.Cells.Item(7,2).Value()
$ws.Cells.Item(7,2).Value = "_10301_Book_A"
$ws.Cells.Item(7,3).Value = "AssetAddress.Book_A_Pickupable"

# --- Row 8: Book B (no assetAddress) ---
$ws.Cells.Item(8,4).Value = $ws.Cells.Item(8,2).Value()
$ws.Cells.Item(8,2).Value = "_10302_Book_B"

# --- Row 9: Book C (no assetAddress) ---
$ws.Cells.Item(9,4).Value = $ws.Cells.Item(9,2).Value()
$ws.Cells.Item(9,2).Value = "_10302_Book_C"

# --- Row 10: drop the old "Coin" row entirely ---
$ws.Cells.Item(10,1).ClearContents()
$ws.Cells.Item(10,2).ClearContents()

# --- Column widths (no longer auto best-fit; explicit widths instead) ---
$ws.Columns.Item(2).ColumnWidth = 16.67
$ws.Columns.Item(3).ColumnWidth = 38.33

# --- Add the new "Notes" sheet after "Entities" ---
$notes = $wb.Worksheets.Add($null, $ws)
$notes.Name = "Notes"
$notes.Range("B2").Value = "title is not being used since it will get localized in the MasterLocalization"

# Copy H1's "Note" cell style (fill/border) onto B2:H2 without disturbing the value we just set
$ws.Range("H1").Copy()
$notes.Range("B2:H2").PasteSpecial(-4122)

# --- Selections: leave "Entities" as the active/selected tab, matching each sheet's own cursor ---
[void]$notes.Range("D7").Select()
[void]$ws.Range("D12").Select()

Write-Output "done"
